# =========================================================================
# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets
#
# This script:
#   1. Inserts a new "Player Info" sheet as the FIRST sheet, with the
#      player's id/name/batting-hand/bowling-style.
#   2. Keeps the existing "ODI Batting" and "ODI Bowling" sheets (now 2nd
#      and 3rd respectively), but:
#        - renames their MATCH_CARD_LINK column to MATCH_CODE
#        - replaces the full howstat.com URL values with just the numeric
#          match code that was at the end of the URL
#        - removes a handful of stray empty cells in column B of
#          "ODI Batting" (rows where the batter did not bat, so there is
#          no INNING_NUMBER)
#   3. Appends a new "ODI Batting Extra" sheet as the LAST sheet, with
#      additional per-match batting detail for the player.
# =========================================================================

$wb = $excel.ActiveWorkbook

$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

# -------------------------------------------------------------------------
# Add "ODI Batting Extra" at the very end first (while sheet order/
# references are still simple), then insert "Player Info" before
# "ODI Batting" so it becomes sheet #1. Doing the "append at end" step
# first keeps the Before/After sheet references unambiguous.
# -------------------------------------------------------------------------
$wsExtra = $wb.Worksheets.Add($null, $wsBowling)
$wsExtra.Name = "ODI Batting Extra"

$wsPlayerInfo = $wb.Worksheets.Add($wsBatting)
$wsPlayerInfo.Name = "Player Info"

# =========================================================================
# 1) Player Info sheet
# =========================================================================
$wsPlayerInfo.Range("A1").Value = "ID"
$wsPlayerInfo.Range("B1").Value = "NAME"
$wsPlayerInfo.Range("C1").Value = "BATTING_HAND"
$wsPlayerInfo.Range("D1").Value = "BOWL_STYLE"

$piHeader = $wsPlayerInfo.Range("A1:D1")
$piHeader.Font.Bold = $true
$piHeader.Borders.LineStyle = 1
$piHeader.HorizontalAlignment = -4108
$piHeader.VerticalAlignment = -4160

$wsPlayerInfo.Range("A2").NumberFormat = "@"
$wsPlayerInfo.Range("A2").Value = "4605"
$wsPlayerInfo.Range("B2").Value = "Rovman Powell"
$wsPlayerInfo.Range("C2").Value = "Right Handed"
$wsPlayerInfo.Range("D2").Value = "Right Arm Medium Fast"

# =========================================================================
# 2) ODI Batting sheet: MATCH_CARD_LINK -> MATCH_CODE, strip URL prefix,
#    and drop stray empty INNING_NUMBER cells.
# =========================================================================
$wsBatting.Range("D1").Value = "MATCH_CODE"

$battingLastRow = 48
for ($r = 2; $r -le $battingLastRow; $r++) {
    $linkCell = $wsBatting.Cells.Item($r, 4)
    $link = $linkCell.Value2
    $code = $link -replace '.*MatchCode=', ''
    $linkCell.NumberFormat = "@"
    $linkCell.Value = $code

    $inningCell = $wsBatting.Cells.Item($r, 2)
    if ($inningCell.Value2 -eq "") {
        $inningCell.ClearContents()
    }
}

# =========================================================================
# 3) ODI Bowling sheet: MATCH_CARD_LINK -> MATCH_CODE, strip URL prefix.
# =========================================================================
$wsBowling.Range("B1").Value = "MATCH_CODE"

$bowlingLastRow = 16
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $linkCell = $wsBowling.Cells.Item($r, 2)
    $link = $linkCell.Value2
    $code = $link -replace '.*MatchCode=', ''
    $linkCell.NumberFormat = "@"
    $linkCell.Value = $code
}

# =========================================================================
# 4) ODI Batting Extra sheet: new per-match batting detail.
# =========================================================================
$wsExtra.Range("A1").Value = "MATCH_CODE"
$wsExtra.Range("B1").Value = "BATTING_POSITION"
$wsExtra.Range("C1").Value = "NUM_4"
$wsExtra.Range("D1").Value = "NUM_6"
$wsExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Range("F1").Value = "MAN_OF_MATCH"

$extraHeader = $wsExtra.Range("A1:F1")
$extraHeader.Font.Bold = $true
$extraHeader.Borders.LineStyle = 1
$extraHeader.HorizontalAlignment = -4108
$extraHeader.VerticalAlignment = -4160

function Set-ExtraRow {
    param(
        [int]$Row,
        [string]$MatchCode,
        [object]$BattingPosition,
        [string]$Num4,
        [string]$Num6,
        [string]$PercentRunsOfTotal,
        [string]$ManOfMatch
    )

    $codeCell = $wsExtra.Cells.Item($Row, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $MatchCode

    $posCell = $wsExtra.Cells.Item($Row, 2)
    $num4Cell = $wsExtra.Cells.Item($Row, 3)
    $num6Cell = $wsExtra.Cells.Item($Row, 4)
    $pctCell = $wsExtra.Cells.Item($Row, 5)

    if ($null -eq $BattingPosition) {
        $posCell.ClearContents()
        $num4Cell.ClearContents()
        $num6Cell.ClearContents()
        $pctCell.ClearContents()
    } else {
        $posCell.Value = $BattingPosition

        $num4Cell.NumberFormat = "@"
        $num4Cell.Value = $Num4

        $num6Cell.NumberFormat = "@"
        $num6Cell.Value = $Num6

        $pctCell.NumberFormat = "@"
        $pctCell.Value = $PercentRunsOfTotal
    }

    $momCell = $wsExtra.Cells.Item($Row, 6)
    $momCell.Value = $ManOfMatch
}

Set-ExtraRow -Row 2  -MatchCode "4216" -BattingPosition 6    -Num4 "1" -Num6 "1" -PercentRunsOfTotal "5.61%"  -ManOfMatch "NO"
Set-ExtraRow -Row 3  -MatchCode "4219" -BattingPosition 6    -Num4 "0" -Num6 "0" -PercentRunsOfTotal "1.41%"  -ManOfMatch "NO"
Set-ExtraRow -Row 4  -MatchCode "4220" -BattingPosition 6    -Num4 "0" -Num6 "0" -PercentRunsOfTotal "0.65%"  -ManOfMatch "NO"
Set-ExtraRow -Row 5  -MatchCode "4221" -BattingPosition 2    -Num4 "1" -Num6 "0" -PercentRunsOfTotal "15.38%" -ManOfMatch "NO"
Set-ExtraRow -Row 6  -MatchCode "4228" -BattingPosition 6    -Num4 "1" -Num6 "0" -PercentRunsOfTotal "7.18%"  -ManOfMatch "NO"
Set-ExtraRow -Row 7  -MatchCode "4229" -BattingPosition 6    -Num4 "0" -Num6 "0" -PercentRunsOfTotal "0.39%"  -ManOfMatch "NO"
Set-ExtraRow -Row 8  -MatchCode "4230" -BattingPosition 6    -Num4 "0" -Num6 "0" -PercentRunsOfTotal "0.51%"  -ManOfMatch "NO"
Set-ExtraRow -Row 9  -MatchCode "4443" -BattingPosition $null                                                 -ManOfMatch "NO"
Set-ExtraRow -Row 10 -MatchCode "4445" -BattingPosition 8    -Num4 "2" -Num6 "1" -PercentRunsOfTotal "27.70%" -ManOfMatch "NO"
Set-ExtraRow -Row 11 -MatchCode "4447" -BattingPosition 6    -Num4 "2" -Num6 "2" -PercentRunsOfTotal "26.55%" -ManOfMatch "NO"
Set-ExtraRow -Row 12 -MatchCode "4586" -BattingPosition $null                                                 -ManOfMatch "NO"
Set-ExtraRow -Row 13 -MatchCode "4590" -BattingPosition $null                                                 -ManOfMatch "NO"
Set-ExtraRow -Row 14 -MatchCode "4592" -BattingPosition 6    -Num4 "1" -Num6 "0" -PercentRunsOfTotal "4.63%"  -ManOfMatch "NO"
Set-ExtraRow -Row 15 -MatchCode "4606" -BattingPosition $null                                                 -ManOfMatch "NO"
Set-ExtraRow -Row 16 -MatchCode "4611" -BattingPosition 6    -Num4 "2" -Num6 "0" -PercentRunsOfTotal "12.04%" -ManOfMatch "NO"
Set-ExtraRow -Row 17 -MatchCode "4616" -BattingPosition 6    -Num4 "0" -Num6 "1" -PercentRunsOfTotal "10.11%" -ManOfMatch "NO"
Set-ExtraRow -Row 18 -MatchCode "4621" -BattingPosition 6    -Num4 "1" -Num6 "0" -PercentRunsOfTotal "1.97%"  -ManOfMatch "NO"
Set-ExtraRow -Row 19 -MatchCode "4623" -BattingPosition $null                                                 -ManOfMatch "NO"
Set-ExtraRow -Row 20 -MatchCode "4727" -BattingPosition $null                                                 -ManOfMatch "NO"
Set-ExtraRow -Row 21 -MatchCode "4731" -BattingPosition 6    -Num4 "0" -Num6 "0" -PercentRunsOfTotal "0.77%"  -ManOfMatch "NO"

Write-Host "Applied ADDITIONAL SCRAPING edit."
